# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-15 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-16 Tuesday", 2) | Out-Null

# Update the division-problem answer table. The table has 20 rows but only
# every 4th row (1, 5, 9, 13, 17) actually holds data across 5 columns.
$t = $d.Tables.Item(1)

$newValues = @{
    "1,1"  = "85÷7=12, 1"
    "1,2"  = "11÷8=1, 3"
    "1,3"  = "26÷9=2, 8"
    "1,4"  = "99÷8=12, 3"
    "1,5"  = "45÷8=5, 5"

    "5,1"  = "66÷2=33, 0"
    "5,2"  = "90÷3=30, 0"
    "5,3"  = "44÷8=5, 4"
    "5,4"  = "57÷8=7, 1"
    "5,5"  = "99÷7=14, 1"

    "9,1"  = "98÷7=14, 0"
    "9,2"  = "19÷3=6, 1"
    "9,3"  = "10÷3=3, 1"
    "9,4"  = "36÷9=4, 0"
    "9,5"  = "17÷6=2, 5"

    "13,1" = "19÷4=4, 3"
    "13,2" = "74÷5=14, 4"
    "13,3" = "34÷8=4, 2"
    "13,4" = "90÷7=12, 6"
    "13,5" = "84÷8=10, 4"

    "17,1" = "59÷4=14, 3"
    "17,2" = "62÷8=7, 6"
    "17,3" = "48÷2=24, 0"
    "17,4" = "56÷8=7, 0"
    "17,5" = "39÷9=4, 3"
}

foreach ($row in @(1, 5, 9, 13, 17)) {
    for ($col = 1; $col -le 5; $col++) {
        $key = "$row,$col"
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $newValues[$key]
    }
}
